# Update "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" worksheets, which contain identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 284
    $ws.Range("F4").Value = 98
    $ws.Range("F5").Value = 873
    $ws.Range("F6").Value = 211
}
